$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Workbook-level view / calc settings (best-effort; some of these are not
#    persisted by the host's OOXML writer, but they are the correct COM verbs)
# ---------------------------------------------------------------------------
$excel.Iteration = $true

# ---------------------------------------------------------------------------
# 2. Add the three new "Paradox ID" -> "IO Name" rows at the bottom of the
#    table (rows 26-28), before the data gets re-sorted by column B.
#    Column A on rows 26 and 28 needs the "Normal 4" text style that's used
#    by other rows of this same flavour (e.g. row 21) - copy that formatting
#    across first, then overwrite the value so the style index is preserved.
# ---------------------------------------------------------------------------
$ws.Range("A21").Copy()
$ws.Range("A26").PasteSpecial(-4122)
$ws.Range("A26").Value = '~RackA\SGr1\RefLevel Measure'
$ws.Range("B26").Value = 'Refrigerant Level `%rackname`'

$ws.Range("B27").Value = 'Comp VFD Value `%rackname`'
$ws.Range("A27").Value = 'Rack B VFD % Measure'

$ws.Range("A21").Copy()
$ws.Range("A28").PasteSpecial(-4122)
$ws.Range("A28").Value = '~A01 POS Freezer #1\Defrost Temp'
$ws.Range("B28").Value = 'System Defrost Temp `%rackname` `%sgname` `%sysname`'

# ---------------------------------------------------------------------------
# 3. Re-sort the whole table (A2:B28) by column B ascending, same as the
#    worksheet's persisted sortState/sortCondition.
# ---------------------------------------------------------------------------
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("B2:B28"))
$ws.Sort.SetRange($ws.Range("A2:B28"))
$ws.Sort.Header = 0
$ws.Sort.Apply()

# ---------------------------------------------------------------------------
# 4. Restore the active cell selection shown in the file.
# ---------------------------------------------------------------------------
$ws.Range("B12").Select()
